# Applies the diff:
#  - Removes columns: discovered_employees, discovered_revenue, discovered_industry, flagged_rpe
#    (originally columns CB:CE, i.e. 80-83)
#  - Inserts two new columns right after "Reason_for_domain_mismatch_mail":
#       WorkPhone_Reason       (blank value)
#       WorkPhone_ColorFlag    (value "False")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the four obsolete columns (CB:CE -> discovered_employees, discovered_revenue,
#    discovered_industry, flagged_rpe).
$ws.Columns("CB:CE").Delete()

# After the deletion, "Reason_for_domain_mismatch_mail" sits at column CD (82) and
# "linkedin_link_found" sits at column CE (83). Insert two blank columns before
# "linkedin_link_found" so the new fields land right after "Reason_for_domain_mismatch_mail".
$ws.Columns("CE:CF").Insert()

# 2) Populate the two new header cells + the data row values.
$ws.Cells.Item(1, 83).Value = "WorkPhone_Reason"
$ws.Cells.Item(1, 84).Value = "WorkPhone_ColorFlag"

$ws.Cells.Item(2, 84).Value = "False"
